# The original workbook reuses two shared-string entries (the Handoff/Handback
# timestamps for the 618f72fa file) across two rows each, on both the zh-cn and
# de-de sheets: row 3 (618f72fa-...) and row 5 (dcff0ae8-...) happen to have had
# identical Correspond Handoff/Handback datetimes, so both rows point at the same
# shared-string entries. The commit bumps those shared strings in place, which
# changes the displayed value everywhere it is referenced - i.e. in both rows.
$wb = $excel.ActiveWorkbook

# zh-cn sheet: rows 3 and 5 shared "2016-03-12 16:16:16" / "2016-03-12 16:16:30"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-12 16:16:57"
$wsZhCn.Range("H3").Value = "2016-03-12 16:17:14"
$wsZhCn.Range("E5").Value = "2016-03-12 16:16:57"
$wsZhCn.Range("H5").Value = "2016-03-12 16:17:14"

# de-de sheet: rows 3 and 5 shared "2016-03-12 16:16:19" / "2016-03-12 16:16:35"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-12 16:17:01"
$wsDeDe.Range("H3").Value = "2016-03-12 16:17:19"
$wsDeDe.Range("E5").Value = "2016-03-12 16:17:01"
$wsDeDe.Range("H5").Value = "2016-03-12 16:17:19"
